{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Version 1.0\" paragraph and the empty paragraph that\n// immediately precedes \"Jack Hulspas\" (the last of the run of blank,\n// colored placeholder paragraphs on the title page).\nlet versionIdx = -1;\nlet jackIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"Version 1.0\") versionIdx = i;\n  if (t === \"Jack Hulspas\") jackIdx = i;\n}\n\nif (versionIdx === -1 || jackIdx === -1) {\n  throw new Error(\"Could not locate the expected title-page paragraphs.\");\n}\n\n// Insert the new date paragraph right after \"Version 1.0\", matching its\n// character formatting (12pt, no special color).\nconst versionPara = items[versionIdx];\nconst dateText = \"09/08/2023 (DD/MM/YYYY)\";\nconst newPara = versionPara.insertParagraph(dateText, \"After\");\nnewPara.font.size = 12;\n\n// Remove the last blank placeholder paragraph right before \"Jack Hulspas\".\nconst blankPara = items[jackIdx - 1];\nblankPara.load(\"text\");\nawait context.sync();\nif (blankPara.text.trim() === \"\") {\n  blankPara.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Insert the new date paragraph right after \"Version 1.0\" ---------------\n$versionRange = $d.Content\n$versionRange.Find.ClearFormatting()\n$foundVersion = $versionRange.Find.Execute(\"Version 1.0\")\nif (-not $foundVersion) {\n    throw \"Could not find the 'Version 1.0' paragraph.\"\n}\n$versionRange.Collapse(0)   # wdCollapseEnd\n$versionRange.InsertParagraphAfter()\n\n$versionPara = $versionRange.Paragraphs(1)\n$newPara = $versionPara.Next()\n$newPara.Range.Text = \"09/08/2023 (DD/MM/YYYY)\"\n\n# --- Remove the blank placeholder paragraph right before \"Jack Hulspas\" ----\n$jackRange = $d.Content\n$jackRange.Find.ClearFormatting()\n$foundJack = $jackRange.Find.Execute(\"Jack Hulspas\")\nif (-not $foundJack) {\n    throw \"Could not find the 'Jack Hulspas' paragraph.\"\n}\n$jackPara = $jackRange.Paragraphs(1)\n$blankPara = $jackPara.Previous()\nif ($blankPara.Range.Text.Trim() -eq \"\") {\n    $blankPara.Range.Delete()\n}\n\n$d.Save()\n"}
